$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.544.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.84%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.693.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.31%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'557.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.23%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'159.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.19%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.86%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -3.79%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.20%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.370"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.84%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'5.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -8.26%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.170.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.38%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'26.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.81%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'63.392.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.40%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -3.68%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.696.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.48%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'12.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.00%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -5.22%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'346.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.73%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -4.31%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.23%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -3.72%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'63.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.75%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.49%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.12%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'8.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.23%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0862"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.58%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.20%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.48%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.44%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'165.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.86%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.05%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -2.43%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -2.92%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -2.62%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -1.11%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'344.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.18%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.951"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -5.30%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'6.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.90%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'38.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.08%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -5.49%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'InjectiveProtocol"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'20.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -4.60%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'EnergySwap"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'20.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -5.46%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.68%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0565"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.70%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.02%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'11.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.02%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'130.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.67%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -3.70%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -4.58%  "
$ws.Range("E51").Style = "Normal"

Write-Output "Applied all crypto list updates"
